# Update countries & provincias Spain
# - Reorders two pairs/blocks of country rows (Barbados/Uganda and
#   Mali ahead of Etiopia..Republica de Yibuti)
# - Refreshes the "Datos actualizados" timestamp
# - Updates the day's case/recovered/death counters for several countries

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Datos actualizados" timestamp (cell A1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 13:20"

# --- Re-order country names ---
# Barbados now sits before Uganda (rows 136-137)
$ws.Cells.Item(136, 1).Value = "Barbados"
$ws.Cells.Item(137, 1).Value = "Uganda"

# Mali moves ahead of Etiopia, shifting the intervening countries down one row
# (rows 143-148): Etiopia, Bermudas, Niger, Congo, Republica de Yibuti, Mali
# becomes: Mali, Etiopia, Bermudas, Niger, Congo, Republica de Yibuti
$ws.Cells.Item(143, 1).Value = "Mali"
$ws.Cells.Item(144, 1).Value = "Etiopia"
$ws.Cells.Item(145, 1).Value = "Bermudas"
$ws.Cells.Item(146, 1).Value = "Niger"
$ws.Cells.Item(147, 1).Value = "Congo"
$ws.Cells.Item(148, 1).Value = "Republica de Yibuti"

# --- Update numeric counters (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) for affected rows ---

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 142793
$ws.Cells.Item(4, 3).Value = 333
$ws.Cells.Item(4, 4).Value = 4562
$ws.Cells.Item(4, 5).Value = 135741
$ws.Cells.Item(4, 6).Value = 2970
$ws.Cells.Item(4, 7).Value = 6
$ws.Cells.Item(4, 8).Value = 2490

# Row 17: Austria
$ws.Cells.Item(17, 2).Value = 9200
$ws.Cells.Item(17, 3).Value = 412
$ws.Cells.Item(17, 4).Value = 636
$ws.Cells.Item(17, 5).Value = 8456

# Row 25: Chequia
$ws.Cells.Item(25, 2).Value = 2866
$ws.Cells.Item(25, 3).Value = 49
$ws.Cells.Item(25, 4).Value = 11
$ws.Cells.Item(25, 5).Value = 2838

# Row 74: Bosnia y Herzegovina
$ws.Cells.Item(74, 2).Value = 354
$ws.Cells.Item(74, 3).Value = 31
$ws.Cells.Item(74, 4).Value = 17
$ws.Cells.Item(74, 5).Value = 329
$ws.Cells.Item(74, 6).Value = 1
$ws.Cells.Item(74, 7).Value = 2
$ws.Cells.Item(74, 8).Value = 8

# Row 93: Oman
$ws.Cells.Item(93, 2).Value = 179
$ws.Cells.Item(93, 3).Value = 12
$ws.Cells.Item(93, 4).Value = 29
$ws.Cells.Item(93, 5).Value = 150

# Row 99: Uzbekistan
$ws.Cells.Item(99, 2).Value = 149
$ws.Cells.Item(99, 3).Value = 5
$ws.Cells.Item(99, 4).Value = 7
$ws.Cells.Item(99, 5).Value = 140

# Row 124: Liechtenstein
$ws.Cells.Item(124, 2).Value = 62
$ws.Cells.Item(124, 3).Value = 6
$ws.Cells.Item(124, 4).Value = 0
$ws.Cells.Item(124, 5).Value = 62

# Row 143: now Mali
$ws.Cells.Item(143, 2).Value = 25
$ws.Cells.Item(143, 3).Value = 7
$ws.Cells.Item(143, 4).Value = 0
$ws.Cells.Item(143, 5).Value = 23
$ws.Cells.Item(143, 6).Value = 0
$ws.Cells.Item(143, 7).Value = 1
$ws.Cells.Item(143, 8).Value = 2

# Row 144: now Etiopia
$ws.Cells.Item(144, 2).Value = 23
$ws.Cells.Item(144, 3).Value = 2
$ws.Cells.Item(144, 4).Value = 1
$ws.Cells.Item(144, 5).Value = 22

# Row 145: now Bermudas
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(145, 4).Value = 2
$ws.Cells.Item(145, 5).Value = 20
$ws.Cells.Item(145, 6).Value = 0
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 0

# Row 146: now Niger
$ws.Cells.Item(146, 2).Value = 22
$ws.Cells.Item(146, 3).Value = 4
$ws.Cells.Item(146, 7).Value = 2
$ws.Cells.Item(146, 8).Value = 3

# Row 147: now Congo
$ws.Cells.Item(147, 2).Value = 19
$ws.Cells.Item(147, 5).Value = 19

# Row 148: now Republica de Yibuti
$ws.Cells.Item(148, 5).Value = 18
$ws.Cells.Item(148, 8).Value = 0
